$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Dinesh Chandimal"
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "LBW"
$ws.Range("E2").Value = " Mark Wood"
$ws.Range("J2").Value = "Jason Roy"
$ws.Range("K2").Value = 27
$ws.Range("L2").Value = 11
$ws.Range("M2").Value = "LBW"
$ws.Range("N2").Value = " Dushmantha Chameera"
# Row 3
$ws.Range("A3").Value = "Pathum Nissanka"
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 10
$ws.Range("D3").Value = "Bowled"
$ws.Range("E3").Value = " Chris Woakes"
$ws.Range("J3").Value = "Jos Buttler"
$ws.Range("K3").Value = 21
$ws.Range("L3").Value = 7
$ws.Range("M3").Value = "Bowled"
$ws.Range("N3").Value = " Chamika Karunarathne"
# Row 4
$ws.Range("A4").Value = "Charith Asalanka"
$ws.Range("B4").Value = 26
$ws.Range("C4").Value = 9
$ws.Range("E4").Value = " Chris Jordan"
$ws.Range("J4").Value = "Dawid Malan"
$ws.Range("K4").Value = 7
$ws.Range("L4").Value = 4
$ws.Range("M4").Value = "NOT OUT"
$ws.Range("N4").Value = " "
# Row 5
$ws.Range("A5").Value = "Dhananjaya de Silva"
$ws.Range("B5").Value = 1
$ws.Range("D5").Value = "Bowled"
$ws.Range("E5").Value = " Chris Woakes"
$ws.Range("J5").Value = "Jonny Bairstow"
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 3
$ws.Range("M5").Value = "Caught"
$ws.Range("N5").Value = " Chamika Karunarathne"
# Row 6
$ws.Range("A6").Value = "Bhanuka Rajapakse"
$ws.Range("B6").Value = 27
$ws.Range("C6").Value = 10
$ws.Range("D6").Value = "Bowled"
$ws.Range("E6").Value = " Adil Rashid"
$ws.Range("J6").Value = "Eoin Morgan(C)"
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = "Bowled"
$ws.Range("N6").Value = " Chamika Karunarathne"
# Row 7
$ws.Range("A7").Value = "Dasun Shanka(C)"
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 3
$ws.Range("D7").Value = "Bowled"
$ws.Range("E7").Value = " Liam Livingstone"
$ws.Range("J7").Value = "Moeen Ali"
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = "LBW"
$ws.Range("N7").Value = " Chamika Karunarathne"
# Row 8
$ws.Range("A8").Value = "Wanindu Hasaranga"
$ws.Range("B8").Value = 38
$ws.Range("C8").Value = 10
$ws.Range("D8").Value = "LBW"
$ws.Range("E8").Value = " Mark Wood"
$ws.Range("J8").Value = "Liam Livingstone"
$ws.Range("K8").Value = 14
$ws.Range("L8").Value = 7
$ws.Range("M8").Value = "Caught"
$ws.Range("N8").Value = " Nuwan Pradeep"
# Row 9
$ws.Range("A9").Value = "Chamika Karunarathne"
$ws.Range("B9").Value = 14
$ws.Range("C9").Value = 4
$ws.Range("D9").Value = "Caught"
$ws.Range("E9").Value = " Chris Jordan"
$ws.Range("J9").Value = "Chris Woakes"
$ws.Range("K9").Value = 10
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = "LBW"
$ws.Range("N9").Value = " Nuwan Pradeep"
# Row 10
$ws.Range("A10").Value = "Dushmantha Chameera"
$ws.Range("B10").Value = 8
$ws.Range("C10").Value = 3
$ws.Range("E10").Value = " Chris Jordan"
$ws.Range("J10").Value = "Chris Jordan"
$ws.Range("K10").Value = 7
$ws.Range("L10").Value = 3
$ws.Range("M10").Value = "Bowled"
$ws.Range("N10").Value = " Maheesh Theekshana"
# Row 11
$ws.Range("A11").Value = "Maheesh Theekshana"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = "Bowled"
$ws.Range("E11").Value = " Chris Woakes"
$ws.Range("J11").Value = "Adil Rashid"
$ws.Range("K11").Value = 8
$ws.Range("L11").Value = 5
$ws.Range("M11").Value = "Bowled"
$ws.Range("N11").Value = " Dushmantha Chameera"
# Row 12
$ws.Range("A12").Value = "Nuwan Pradeep"
$ws.Range("B12").Value = 26
$ws.Range("C12").Value = 9
$ws.Range("D12").Value = "NOT OUT"
$ws.Range("E12").Value = " "
$ws.Range("J12").Value = "Mark Wood"
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = "LBW"
$ws.Range("N12").Value = " Dushmantha Chameera"
# Row 16
$ws.Range("A16").Value = 174
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "10.2"
$ws.Range("D16").Value = 62
$ws.Range("J16").Value = 101
$ws.Range("K16").Value = 10
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "7.4"
$ws.Range("M16").Value = 46
# Row 21
$ws.Range("A21").Value = "Adil Rashid"
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2.0"
$ws.Range("C21").Value = 34
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = 17
$ws.Range("J21").Value = "Chamika Karunarathne"
$ws.Range("L21").Value = 7
$ws.Range("M21").Value = 4
$ws.Range("N21").Value = 7
# Row 22
$ws.Range("A22").Value = "Chris Jordan"
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "2.0"
$ws.Range("C22").Value = 30
$ws.Range("D22").Value = 3
$ws.Range("J22").Value = "Wanindu Hasaranga"
$ws.Range("L22").Value = 15
$ws.Range("M22").Value = 0
$ws.Range("N22").Value = 15
# Row 23
$ws.Range("A23").Value = "Chris Woakes"
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "2.0"
$ws.Range("C23").Value = 32
$ws.Range("D23").Value = 3
$ws.Range("E23").Value = 16
$ws.Range("J23").Value = "Nuwan Pradeep"
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = "2.0"
$ws.Range("L23").Value = 31
$ws.Range("M23").Value = 2
$ws.Range("N23").Value = 15.5
# Row 24
$ws.Range("A24").Value = "Liam Livingstone"
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "2.0"
$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 15.5
$ws.Range("J24").Value = "Maheesh Theekshana"
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = "2.0"
$ws.Range("L24").Value = 34
$ws.Range("M24").Value = 1
$ws.Range("N24").Value = 17
# Row 25
$ws.Range("A25").Value = "Mark Wood"
$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "2.2"
$ws.Range("C25").Value = 47
$ws.Range("D25").Value = 2
$ws.Range("E25").Value = 21.36
$ws.Range("J25").Value = "Dushmantha Chameera"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = "1.4"
$ws.Range("L25").Value = 14
$ws.Range("M25").Value = 3
$ws.Range("N25").Value = 10

Write-Output "done"